$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.160.96'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.422.68'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '553.52'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.57'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.41%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.588'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.71%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.69'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.09%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.28%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '25.08'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.855.11'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.082.04'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.440.44'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.32'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.50'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '328.76'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.74'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.54%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.70'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.72'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.10%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0775'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.21%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '170.32'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.12'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -3.53%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.61%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.406'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -4.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.59'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.08%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.26%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.22'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '334.16'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.67%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.81%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '38.85'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '146.48'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.66%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.49%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.09'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.41%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0967'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0517'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.96%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.579'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0223'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '11.04'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.58'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.90%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.05%  '
